{"js": "// Apply targeted text replacements (date header + multiplication table cells)\n// using Office.js search/replace so existing run formatting (fonts, size,\n// paragraph alignment, etc.) is preserved untouched.\nconst replacements = [\n  [\"2025-05-25 Sunday\", \"2025-05-26 Monday\"],\n  [\"333\u00d76=1998\", \"479\u00d75=2395\"],\n  [\"576\u00d75=2880\", \"913\u00d79=8217\"],\n  [\"568\u00d74=2272\", \"590\u00d78=4720\"],\n  [\"647\u00d78=5176\", \"781\u00d73=2343\"],\n  [\"173\u00d73=519\", \"812\u00d77=5684\"],\n  [\"362\u00d76=2172\", \"215\u00d73=645\"],\n  [\"537\u00d74=2148\", \"887\u00d73=2661\"],\n  [\"524\u00d75=2620\", \"809\u00d72=1618\"],\n  [\"704\u00d74=2816\", \"597\u00d73=1791\"],\n  [\"967\u00d73=2901\", \"559\u00d76=3354\"],\n  [\"300\u00d72=600\", \"806\u00d78=6448\"],\n  [\"302\u00d75=1510\", \"154\u00d79=1386\"],\n  [\"992\u00d74=3968\", \"358\u00d76=2148\"],\n  [\"234\u00d72=468\", \"638\u00d76=3828\"],\n  [\"877\u00d79=7893\", \"299\u00d73=897\"],\n  [\"960\u00d76=5760\", \"711\u00d79=6399\"],\n  [\"403\u00d77=2821\", \"934\u00d79=8406\"],\n  [\"252\u00d78=2016\", \"564\u00d76=3384\"],\n  [\"172\u00d79=1548\", \"203\u00d72=406\"],\n  [\"703\u00d72=1406\", \"574\u00d74=2296\"],\n  [\"321\u00d79=2889\", \"819\u00d72=1638\"],\n  [\"156\u00d76=936\", \"685\u00d77=4795\"],\n  [\"967\u00d72=1934\", \"225\u00d74=900\"],\n  [\"783\u00d75=3915\", \"839\u00d77=5873\"],\n  [\"104\u00d72=208\", \"906\u00d73=2718\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply targeted text replacements (date header + multiplication table cells)\n# via Word COM Find/Replace, preserving existing run formatting since only\n# the matched text range is replaced in place.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-25 Sunday\", \"2025-05-26 Monday\"),\n    @(\"333\u00d76=1998\", \"479\u00d75=2395\"),\n    @(\"576\u00d75=2880\", \"913\u00d79=8217\"),\n    @(\"568\u00d74=2272\", \"590\u00d78=4720\"),\n    @(\"647\u00d78=5176\", \"781\u00d73=2343\"),\n    @(\"173\u00d73=519\", \"812\u00d77=5684\"),\n    @(\"362\u00d76=2172\", \"215\u00d73=645\"),\n    @(\"537\u00d74=2148\", \"887\u00d73=2661\"),\n    @(\"524\u00d75=2620\", \"809\u00d72=1618\"),\n    @(\"704\u00d74=2816\", \"597\u00d73=1791\"),\n    @(\"967\u00d73=2901\", \"559\u00d76=3354\"),\n    @(\"300\u00d72=600\", \"806\u00d78=6448\"),\n    @(\"302\u00d75=1510\", \"154\u00d79=1386\"),\n    @(\"992\u00d74=3968\", \"358\u00d76=2148\"),\n    @(\"234\u00d72=468\", \"638\u00d76=3828\"),\n    @(\"877\u00d79=7893\", \"299\u00d73=897\"),\n    @(\"960\u00d76=5760\", \"711\u00d79=6399\"),\n    @(\"403\u00d77=2821\", \"934\u00d79=8406\"),\n    @(\"252\u00d78=2016\", \"564\u00d76=3384\"),\n    @(\"172\u00d79=1548\", \"203\u00d72=406\"),\n    @(\"703\u00d72=1406\", \"574\u00d74=2296\"),\n    @(\"321\u00d79=2889\", \"819\u00d72=1638\"),\n    @(\"156\u00d76=936\", \"685\u00d77=4795\"),\n    @(\"967\u00d72=1934\", \"225\u00d74=900\"),\n    @(\"783\u00d75=3915\", \"839\u00d77=5873\"),\n    @(\"104\u00d72=208\", \"906\u00d73=2718\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $r = $d.Content\n    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
